$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These D-column cells get new values that look numeric (e.g. "482.37").
# The source workbook stores every cell in this sheet as text (inline
# strings), so force a text number format before assigning the value to
# stop Excel from auto-converting the text into a real number.
$numericLookingRefs = @("D5", "D6", "D10", "D11", "D13", "D15", "D18", "D21", "D22", "D23", "D24", "D25", "D27", "D28", "D30", "D31", "D32", "D35", "D36", "D39", "D42", "D43", "D44", "D46", "D49")
foreach ($ref in $numericLookingRefs) {
    $ws.Range($ref).NumberFormat = "@"
}

# Apply the updated coin data (price/volume refresh, plus a couple of rows
# that got reordered and a name/link swap) per the latest scrape.
$ws.Range("D2").Value = '68.269.40'
$ws.Range("E2").Value = '  +1.86%  '
$ws.Range("D3").Value = '3.893.37'
$ws.Range("E3").Value = '  +0.80%  '
$ws.Range("E4").Value = '  +0.23%  '
$ws.Range("D5").Value = '482.37'
$ws.Range("E5").Value = '  +2.06%  '
$ws.Range("D6").Value = '145.27'
$ws.Range("E6").Value = '  +0.31%  '
$ws.Range("E7").Value = '  -1.37%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("E9").Value = '  -2.87%  '
$ws.Range("D10").Value = '0.165'
$ws.Range("E10").Value = '  +6.61%  '
$ws.Range("D11").Value = '0.0000353'
$ws.Range("E11").Value = '  +13.55%  '
$ws.Range("E12").Value = '  -1.67%  '
$ws.Range("D13").Value = '10.64'
$ws.Range("E13").Value = '  +2.19%  '
$ws.Range("D14").Value = '4.525.55'
$ws.Range("E14").Value = '  +1.02%  '
$ws.Range("D15").Value = '14.63'
$ws.Range("E15").Value = '  -1.30%  '
$ws.Range("D16").Value = '3.889.00'
$ws.Range("E16").Value = '  +0.54%  '
$ws.Range("E17").Value = '  -0.33%  '
$ws.Range("D18").Value = '19.76'
$ws.Range("E18").Value = '  -1.76%  '
$ws.Range("E19").Value = '  -3.08%  '
$ws.Range("D20").Value = '68.362.31'
$ws.Range("E20").Value = '  +1.65%  '
$ws.Range("D21").Value = '435.34'
$ws.Range("E21").Value = '  +0.72%  '
$ws.Range("B22").Value = 'ImmutableX'
$ws.Range("C22").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D22").Value = '3.39'
$ws.Range("E22").Value = '  +1.80%  '
$ws.Range("B23").Value = 'InternetComputer(DFINITY)'
$ws.Range("C23").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D23").Value = '14.72'
$ws.Range("E23").Value = '  -1.56%  '
$ws.Range("D24").Value = '87.97'
$ws.Range("E24").Value = '  -0.63%  '
$ws.Range("D25").Value = '11.52'
$ws.Range("E25").Value = '  +15.85%  '
$ws.Range("E26").Value = '  -0.10%  '
$ws.Range("D27").Value = '10.48'
$ws.Range("E27").Value = '  +4.74%  '
$ws.Range("D28").Value = '38.02'
$ws.Range("E28").Value = '  -0.03%  '
$ws.Range("E29").Value = '  +4.57%  '
$ws.Range("D30").Value = '702.88'
$ws.Range("E30").Value = '  -3.80%  '
$ws.Range("B31").Value = 'Cosmos'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D31").Value = '13.36'
$ws.Range("E31").Value = '  -3.85%  '
$ws.Range("B32").Value = 'Hedera'
$ws.Range("C32").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D32").Value = '0.131'
$ws.Range("E32").Value = '  -2.76%  '
$ws.Range("E33").Value = '  +2.77%  '
$ws.Range("D34").Value = '0.0₃0914'
$ws.Range("E34").Value = '  +35.49%  '
$ws.Range("D35").Value = '41.68'
$ws.Range("E35").Value = '  -3.76%  '
$ws.Range("D36").Value = '59.45'
$ws.Range("E36").Value = '  +1.94%  '
$ws.Range("E37").Value = '  +3.68%  '
$ws.Range("E38").Value = '  -6.41%  '
$ws.Range("D39").Value = '0.999'
$ws.Range("E39").Value = '  -0.04%  '
$ws.Range("E40").Value = '  -1.91%  '
$ws.Range("E41").Value = '  +10.04%  '
$ws.Range("D42").Value = '3.04'
$ws.Range("E42").Value = '  +4.24%  '
$ws.Range("D43").Value = '2.75'
$ws.Range("E43").Value = '  +7.29%  '
$ws.Range("D44").Value = '0.343'
$ws.Range("E44").Value = '  -1.22%  '
$ws.Range("E45").Value = '  -0.17%  '
$ws.Range("D46").Value = '1.00'
$ws.Range("E46").Value = '  +0.03%  '
$ws.Range("E47").Value = '  -1.78%  '
$ws.Range("E48").Value = '  -1.10%  '
$ws.Range("D49").Value = '146.18'
$ws.Range("E49").Value = '  +1.82%  '
$ws.Range("E50").Value = '  -2.47%  '
$ws.Range("E51").Value = '  -2.01%  '

# Restore the default (Normal) style on the cells we touched above so no
# stray text-number-format is left applied to them.
foreach ($ref in $numericLookingRefs) {
    $ws.Range($ref).Style = "Normal"
}

Write-Output "Applied cryptos list update"
